$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rewrite the whole used range so that the shared-strings pool comes out
# with "Mr"/"Mukul Sangma"/... kept in their original relative order and
# the two brand-new headers ("Constituency"/"Address") registered right
# after "PartyImageUrl" (before "Mr" is ever touched again), followed by
# the three brand-new location values appended in call order.
$ws.UsedRange.ClearContents()

# --- Row 1: header row ---
$ws.Cells.Item(1, 1).Value = "No"
$ws.Cells.Item(1, 2).Value = "Salutation"
$ws.Cells.Item(1, 3).Value = "FullName"
$ws.Cells.Item(1, 4).Value = "Brief"
$ws.Cells.Item(1, 5).Value = "Education"
$ws.Cells.Item(1, 6).Value = "CurrentDesignation"
$ws.Cells.Item(1, 7).Value = "Born"
$ws.Cells.Item(1, 8).Value = "Parents"
$ws.Cells.Item(1, 9).Value = "Spouse"
$ws.Cells.Item(1, 10).Value = "OfficialWebsite"
$ws.Cells.Item(1, 11).Value = "PoliticalParty"
$ws.Cells.Item(1, 12).Value = "WikipediaUrl"
$ws.Cells.Item(1, 13).Value = "FacebookUrl"
$ws.Cells.Item(1, 14).Value = "InstagramUrl"
$ws.Cells.Item(1, 15).Value = "GooglePlusUrl"
$ws.Cells.Item(1, 16).Value = "LinkedInUrl"
$ws.Cells.Item(1, 17).Value = "TwitterUrl"
$ws.Cells.Item(1, 18).Value = "YoutubeUrl"
$ws.Cells.Item(1, 19).Value = "SpeechUrl"
$ws.Cells.Item(1, 20).Value = "ProfileUrl"
$ws.Cells.Item(1, 21).Value = "PartyImageUrl"
$ws.Cells.Item(1, 22).Value = "Constituency"
$ws.Cells.Item(1, 23).Value = "Address"

# New "Constituency" header cell gets its own font/style (fontId 15 - Arial 10,
# no charset - distinct from the base font).
$ws.Cells.Item(1, 22).Font.Name = "Arial"
$ws.Cells.Item(1, 22).Font.Size = 10

# --- Row 2: Mukul Sangma, record #101 ---
$ws.Cells.Item(2, 1).Value = 101
$ws.Cells.Item(2, 2).Value = "Mr"
$ws.Cells.Item(2, 3).Value = "Mukul Sangma"
$ws.Cells.Item(2, 4).Value = "Mukul M. Sangma is a politician from the Indian National Congress and is a former Chief Minister of Meghalaya."
$ws.Cells.Item(2, 5).Value = "Medicine from the Regional Institute of Medical Sciences"
$ws.Cells.Item(2, 6).Value = "Former Chief minister of Meghalaya"
$ws.Cells.Item(2, 7).Value = "20 April 1965 (age 52)"
$ws.Cells.Item(2, 8).Value = "Roshanara Begum, Binoy Bhushan M Marak"
$ws.Cells.Item(2, 9).Value = "Dikanchi D. Shira"
$ws.Cells.Item(2, 12).Value = "https://en.wikipedia.org/wiki/Mukul_Sangma"
$ws.Cells.Item(2, 13).Value = "http://www.facebook.com/mukulsangma"
$ws.Cells.Item(2, 17).Value = "http://www.twitter.com/mukulsangma"
$ws.Cells.Item(2, 20).Value = "http://localhost:1337/assets/images/MukulSangma_101.jpg"
$ws.Cells.Item(2, 21).Value = "http://localhost:1337/assets/images/IndianNationalCongress.jpg"
$ws.Cells.Item(2, 22).Value = "Baner"
$ws.Cells.Item(2, 23).Value = "Pune"

# --- Row 3: Mukul Sangma, record #102 ---
$ws.Cells.Item(3, 1).Value = 102
$ws.Cells.Item(3, 2).Value = "Mr"
$ws.Cells.Item(3, 3).Value = "Mukul Sangma"
$ws.Cells.Item(3, 4).Value = "Mukul M. Sangma is a politician from the Indian National Congress and is a former Chief Minister of Meghalaya."
$ws.Cells.Item(3, 5).Value = "Medicine from the Regional Institute of Medical Sciences"
$ws.Cells.Item(3, 6).Value = "Former Chief minister of Meghalaya"
$ws.Cells.Item(3, 7).Value = "20 April 1965 (age 52)"
$ws.Cells.Item(3, 8).Value = "Roshanara Begum, Binoy Bhushan M Marak"
$ws.Cells.Item(3, 9).Value = "Dikanchi D. Shira"
$ws.Cells.Item(3, 12).Value = "https://en.wikipedia.org/wiki/Mukul_Sangma"
$ws.Cells.Item(3, 13).Value = "http://www.facebook.com/mukulsangma"
$ws.Cells.Item(3, 17).Value = "http://www.twitter.com/mukulsangma"
$ws.Cells.Item(3, 20).Value = "http://localhost:1337/assets/images/MukulSangma_101.jpg"
$ws.Cells.Item(3, 21).Value = "http://localhost:1337/assets/images/IndianNationalCongress.jpg"
$ws.Cells.Item(3, 22).Value = "Pimpri Chinchwad"
$ws.Cells.Item(3, 23).Value = "Pune"

# Column V (22) is widened to fit "Constituency"/"Baner"/"Pimpri Chinchwad".
$ws.Columns.Item(22).ColumnWidth = 16.02

# Scroll/selection state matches the newly-added columns being in view.
$ws.Range("V3").Select()
